$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: resultado = "Acierto", profit = 1.1
$ws.Range("G10").Value = "Acierto"
$ws.Range("H10").Value = 1.1

# Row 11: resultado = "Fallo", profit = -1
$ws.Range("G11").Value = "Fallo"
$ws.Range("H11").Value = -1
